$wb = $excel.ActiveWorkbook

# --- settings sheet: duplicate the settings block into columns K:R for the
#     new form version, and add a "showFooter" row ---
$ws = $wb.Worksheets.Item("settings")

# Header row (mirrors A1:H1 -> K1:R1)
$ws.Range("K1").Value = "setting_name"
$ws.Range("L1").Value = "value"
$ws.Range("M1").Value = "display.title.text"
$ws.Range("N1").Value = "display.title.text.pt"
$ws.Range("O1").Value = "display.title.text.sw"
$ws.Range("P1").Value = "display.locale.text"
$ws.Range("Q1").Value = "display.locale.text.pt"
$ws.Range("R1").Value = "display.locale.text.sw"

# form_id (unchanged across versions)
$ws.Range("K2").Value = "form_id"
$ws.Range("L2").Value = "hh_geo_location"

# form_version (new version number)
$ws.Range("K3").Value = "form_version"
$ws.Range("L3").Value = 20210221001

# table_id (unchanged across versions)
$ws.Range("K4").Value = "table_id"
$ws.Range("L4").Value = "hh_geo_location"

# survey title (new title)
$ws.Range("K5").Value = "survey"
$ws.Range("M5").Value = "Household Geolocation"
$ws.Range("N5").Value = "Household Geolocation"
$ws.Range("O5").Value = "Household Geolocation"

# locale: default / English
$ws.Range("K6").Value = "default"
$ws.Range("P6").Value = "English"
$ws.Range("Q6").Value = "English"
$ws.Range("R6").Value = "English"

# locale: pt / Português
$ws.Range("K7").Value = "pt"
$ws.Range("P7").Value = "Português"
$ws.Range("Q7").Value = "Português"
$ws.Range("R7").Value = "Português"

# locale: sw / Kiswahili
$ws.Range("K8").Value = "sw"
$ws.Range("P8").Value = "Kiswahili"
$ws.Range("Q8").Value = "Kiswahili"
$ws.Range("R8").Value = "Kiswahili"

# new setting: showFooter = 1 (bottom of screen back/next)
$ws.Range("K9").Value = "showFooter"
$ws.Range("L9").Value = 1

# Make "settings" the active sheet/selection, matching the saved view state.
$ws.Range("L9").Select()
$wb.Worksheets.Item("settings").Activate()
